$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.95578266666667
$ws.Range("H2").Value = 59.867348
$ws.Range("I2").Value = 0.0117373419656925
$ws.Range("J2").Value = 0.0117373419656925
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.297291666666667
$ws.Range("N2").Value = 3.891875
$ws.Range("O2").Value = 0.2261711035865543
$ws.Range("P2").Value = 0.2261711035865543
$ws.Range("Q2").Value = 25.88847055527778
$ws.Range("R2").Value = 232.9962349975
$ws.Range("S2").Value = 0.002654647585553448
$ws.Range("T2").Value = 0.002654647585553449

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.95578266666667
$ws.Range("H3").Value = 59.867348
$ws.Range("I3").Value = 0.0117373419656925
$ws.Range("J3").Value = 0.0117373419656925
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.532647
$ws.Range("N3").Value = 7.597941
$ws.Range("O3").Value = 0.4415441659754047
$ws.Range("P3").Value = 0.4415441659754047
$ws.Range("Q3").Value = 50.54095310338533
$ws.Range("R3").Value = 454.8685779304679
$ws.Range("S3").Value = 0.00518255486900981
$ws.Range("T3").Value = 0.00518255486900981

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.95578266666667
$ws.Range("H4").Value = 59.867348
$ws.Range("I4").Value = 0.0117373419656925
$ws.Range("J4").Value = 0.0117373419656925
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.5649363333333334
$ws.Range("N4").Value = 1.694809
$ws.Range("O4").Value = 0.09849155532960965
$ws.Range("P4").Value = 0.09849155532960965
$ws.Range("Q4").Value = 11.27374668850356
$ws.Range("R4").Value = 101.463720196532
$ws.Range("S4").Value = 0.001156029065636552
$ws.Range("T4").Value = 0.001156029065636552

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.95578266666667
$ws.Range("H5").Value = 59.867348
$ws.Range("I5").Value = 0.0117373419656925
$ws.Range("J5").Value = 0.0117373419656925
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.341011
$ws.Range("N5").Value = 4.023033
$ws.Range("O5").Value = 0.2337931751084314
$ws.Range("P5").Value = 0.2337931751084314
$ws.Range("Q5").Value = 26.76092406960933
$ws.Range("R5").Value = 240.848316626484
$ws.Range("S5").Value = 0.002744110445492686
$ws.Range("T5").Value = 0.002744110445492686

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1637.343343333333
$ws.Range("H6").Value = 4912.03003
$ws.Range("I6").Value = 0.9630320723052701
$ws.Range("J6").Value = 0.9630320723052702
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.297291666666667
$ws.Range("N6").Value = 3.891875
$ws.Range("O6").Value = 0.2261711035865543
$ws.Range("P6").Value = 0.2261711035865543
$ws.Range("Q6").Value = 2124.111874778472
$ws.Range("R6").Value = 19117.00687300625
$ws.Range("S6").Value = 0.2178100265825293
$ws.Range("T6").Value = 0.2178100265825293

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1637.343343333333
$ws.Range("H7").Value = 4912.03003
$ws.Range("I7").Value = 0.9630320723052701
$ws.Range("J7").Value = 0.9630320723052702
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.532647
$ws.Range("N7").Value = 7.597941
$ws.Range("O7").Value = 0.4415441659754047
$ws.Range("P7").Value = 0.4415441659754047
$ws.Range("Q7").Value = 4146.812706463136
$ws.Range("R7").Value = 37321.31435816823
$ws.Range("S7").Value = 0.4252211931735961
$ws.Range("T7").Value = 0.4252211931735961

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1637.343343333333
$ws.Range("H8").Value = 4912.03003
$ws.Range("I8").Value = 0.9630320723052701
$ws.Range("J8").Value = 0.9630320723052702
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.5649363333333334
$ws.Range("N8").Value = 1.694809
$ws.Range("O8").Value = 0.09849155532960965
$ws.Range("P8").Value = 0.09849155532960965
$ws.Range("Q8").Value = 924.9947447904744
$ws.Range("R8").Value = 8324.95270311427
$ws.Range("S8").Value = 0.09485052663364316
$ws.Range("T8").Value = 0.09485052663364317

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1637.343343333333
$ws.Range("H9").Value = 4912.03003
$ws.Range("I9").Value = 0.9630320723052701
$ws.Range("J9").Value = 0.9630320723052702
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.341011
$ws.Range("N9").Value = 4.023033
$ws.Range("O9").Value = 0.2337931751084314
$ws.Range("P9").Value = 0.2337931751084314
$ws.Range("Q9").Value = 2195.695434186776
$ws.Range("R9").Value = 19761.25890768099
$ws.Range("S9").Value = 0.2251503259155015
$ws.Range("T9").Value = 0.2251503259155016

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 17.50081933333334
$ws.Range("H10").Value = 52.502458
$ws.Range("I10").Value = 0.01029341242216722
$ws.Range("J10").Value = 0.01029341242216722
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.297291666666667
$ws.Range("N10").Value = 3.891875
$ws.Range("O10").Value = 0.2261711035865543
$ws.Range("P10").Value = 0.2261711035865543
$ws.Range("Q10").Value = 22.70366708097222
$ws.Range("R10").Value = 204.33300372875
$ws.Range("S10").Value = 0.002328072447193106
$ws.Range("T10").Value = 0.002328072447193107

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 17.50081933333334
$ws.Range("H11").Value = 52.502458
$ws.Range("I11").Value = 0.01029341242216722
$ws.Range("J11").Value = 0.01029341242216722
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.532647
$ws.Range("N11").Value = 7.597941
$ws.Range("O11").Value = 0.4415441659754047
$ws.Range("P11").Value = 0.4415441659754047
$ws.Range("Q11").Value = 44.32339758210867
$ws.Range("R11").Value = 398.910578238978
$ws.Range("S11").Value = 0.004544996202986694
$ws.Range("T11").Value = 0.004544996202986695

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 17.50081933333334
$ws.Range("H12").Value = 52.502458
$ws.Range("I12").Value = 0.01029341242216722
$ws.Range("J12").Value = 0.01029341242216722
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.5649363333333334
$ws.Range("N12").Value = 1.694809
$ws.Range("O12").Value = 0.09849155532960965
$ws.Range("P12").Value = 0.09849155532960965
$ws.Range("Q12").Value = 9.886848704502446
$ws.Range("R12").Value = 88.98163834052201
$ws.Range("S12").Value = 0.001013814199108374
$ws.Range("T12").Value = 0.001013814199108374

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 17.50081933333334
$ws.Range("H13").Value = 52.502458
$ws.Range("I13").Value = 0.01029341242216722
$ws.Range("J13").Value = 0.01029341242216722
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.341011
$ws.Range("N13").Value = 4.023033
$ws.Range("O13").Value = 0.2337931751084314
$ws.Range("P13").Value = 0.2337931751084314
$ws.Range("Q13").Value = 23.46879123501267
$ws.Range("R13").Value = 211.219121115114
$ws.Range("S13").Value = 0.002406529572879043
$ws.Range("T13").Value = 0.002406529572879044

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 25.39612333333333
$ws.Range("H14").Value = 76.18836999999999
$ws.Range("I14").Value = 0.01493717330687017
$ws.Range("J14").Value = 0.01493717330687017
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.297291666666667
$ws.Range("N14").Value = 3.891875
$ws.Range("O14").Value = 0.2261711035865543
$ws.Range("P14").Value = 0.2261711035865543
$ws.Range("Q14").Value = 32.94617916597222
$ws.Range("R14").Value = 296.5156124937499
$ws.Range("S14").Value = 0.003378356971278446
$ws.Range("T14").Value = 0.003378356971278447

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 25.39612333333333
$ws.Range("H15").Value = 76.18836999999999
$ws.Range("I15").Value = 0.01493717330687017
$ws.Range("J15").Value = 0.01493717330687017
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.532647
$ws.Range("N15").Value = 7.597941
$ws.Range("O15").Value = 0.4415441659754047
$ws.Range("P15").Value = 0.4415441659754047
$ws.Range("Q15").Value = 64.31941557179665
$ws.Range("R15").Value = 578.8747401461699
$ws.Range("S15").Value = 0.006595421729812066
$ws.Range("T15").Value = 0.006595421729812066

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 25.39612333333333
$ws.Range("H16").Value = 76.18836999999999
$ws.Range("I16").Value = 0.01493717330687017
$ws.Range("J16").Value = 0.01493717330687017
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.5649363333333334
$ws.Range("N16").Value = 1.694809
$ws.Range("O16").Value = 0.09849155532960965
$ws.Range("P16").Value = 0.09849155532960965
$ws.Range("Q16").Value = 14.34719279681445
$ws.Range("R16").Value = 129.12473517133
$ws.Range("S16").Value = 0.001471185431221572
$ws.Range("T16").Value = 0.001471185431221572

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 25.39612333333333
$ws.Range("H17").Value = 76.18836999999999
$ws.Range("I17").Value = 0.01493717330687017
$ws.Range("J17").Value = 0.01493717330687017
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.341011
$ws.Range("N17").Value = 4.023033
$ws.Range("O17").Value = 0.2337931751084314
$ws.Range("P17").Value = 0.2337931751084314
$ws.Range("Q17").Value = 34.05648074735667
$ws.Range("R17").Value = 306.5083267262099
$ws.Range("S17").Value = 0.003492209174558084
$ws.Range("T17").Value = 0.003492209174558085

